$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: tag the last "27 deg" sample as end of sample ---
$ws.Range("F12").Value = "end of sample"

# --- Row 13: this used to be the (empty-B) "end of sample" row; it is now a
#     real TA reading at 29 deg, tagged "With Junk" ---
$ws.Range("A13").Value = 43188
$ws.Range("B13").Value = 2212.6258499758801
$ws.Range("F13").Value = "With Junk"

# --- Row 14: new TA sample at 29 deg, tagged "Without Junk" ---
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").Value = 43188
$ws.Range("B14").Value = 2223.8240466401598
$ws.Range("C14").Value = 2207.0300000000002
$ws.Range("D14").FormulaR1C1 = "=100*(RC[-2]-RC[-1])/RC[-1]"
$ws.Range("E14").Value = 169
$ws.Range("F14").Value = "Without Junk"

# --- Row 15: new TA sample at 29 deg, tagged both "With Junk" (F) and
#     "end of sample" (G) ---
$ws.Range("A13").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = 43188
$ws.Range("B15").Value = 2225.7790538463801
$ws.Range("F15").Value = "With Junk"
$ws.Range("G15").Value = "end of sample"

$excel.CutCopyMode = 0

$ws.Range("B16").Select()
